$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tiny floating-point re-write of the existing A5 timestamp (matches source diff).
$ws.Range("A5").Value = 44365.96280899306

# New sale row.
$ws.Range("A6").Value = 44365.96721112184
$ws.Range("A6").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("B6").Value = "Heinecken"
$ws.Range("C6").Value = 13
